$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.428.01"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.916.44"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.03"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4072"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08227"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.40"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.922.72"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.075"
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.246"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.41"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06821"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.65"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "29.453.32"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.656"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.78"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D25").Value = "2.170.07"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.664"
$ws.Range("E26").Value = "  +9.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.71"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.018"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09598"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.692"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.552"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.373"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02285"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("E41").Value = "  +6.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1848"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.408"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.279"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07596"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5594"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.953"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.37"
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.425"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.35"
$ws.Range("E51").Value = "  +0.62%  "

# Row 39/40 swap: FraxShare moves above TheSandbox
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.075"
$ws.Range("E39").Value = "  +2.74%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5981"
$ws.Range("E40").Value = "  +2.08%  "
